$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "is_locked_lbl" (D) and "is_enabled_lbl" (E) header columns.
# This shifts the remaining "order_by"/"rem" columns (F,G) left into D,E.
$ws.Range("D1:E1").EntireColumn.Delete()
